$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Liver injury has been described with COVID‐19, and early reports suggested 2%‐11% of patients had chronic liver disease (CLD).`n In this multicentre retrospective study, we evaluated hospitalized adults with laboratory‐confirmed COVID‐19 and the impact of CLD on relevant clinical outcomes.`n Of 363 patients included, 19% had CLD, including 15.2% with NAFLD.`n Patients with CLD had longer length of stay.`n After controlling for age, gender, obesity, cardiac diseases, hypertension, hyperlipidaemia, diabetes and pulmonary disorders, CLD and NAFLD were independently associated with ICU admission ([aOR 1.77, 95% CI 1.03‐3.04] and [aOR 2.30, 95% CI 1.27‐4.17]) and mechanical ventilation ([aOR 2.08, 95% CI 1.20‐3.60] and [aOR 2.15, 95% CI 1.18‐3.91]).`n Presence of cirrhosis was an independent predictor of mortality (aOR 12.5, 95% CI 2.16‐72.5).`n Overall, nearly one‐fifth of hospitalized COVID‐19 patients had CLD, which was associated with more critical illness.`n Future studies are needed to identify interventions to improve clinical outcomes.`n"
$ws.Range("E2").Value = "[Nikroo%Hashemi%NULL%0, Kathleen%Viveiros%NULL%2, Walker D.%Redd%NULL%1, Joyce C.%Zhou%NULL%1, Thomas R.%McCarty%NULL%1, Ahmad N.%Bazarbashi%NULL%1, Kelly E.%Hathorn%NULL%1, Danny%Wong%NULL%2, Cheikh%Njie%NULL%2, Lin%Shen%NULL%2, Walter W.%Chan%wwchan@bwh.harvard.edu%1]"
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = "John Wiley and Sons Inc."
$ws.Range("D3").Value = "Background`nid=`"Par1`">Coronavirus disease-2019 (COVID-19) is a global pandemic.`n`n Obesity has been associated with increased disease severity in COVID-19, and obesity is strongly associated with hepatic steatosis (HS).`n`n However, how HS alters the natural history of COVID-19 is not well characterized, especially in Western populations.`n`n`nAims`nid=`"Par2`">To characterize the impact of HS on disease severity and liver injury in COVID-19.`nMethods`nid=`"Par3`">We examined the association between HS and disease severity in a single-center cohort study of hospitalized COVID-19 patients at Michigan Medicine.`n`n HS was defined by either hepatic steatosis index &gt; 36 (for Asians) or &gt; 39 (for non-Asians) or liver imaging demonstrating steatosis &gt; 30 days before onset of COVID-19. The primary predictor was HS.`n`n The primary outcomes were severity of cardiopulmonary disease, transaminitis, jaundice, and portal hypertensive complications.`n`n`nResults`nid=`"Par4`">In a cohort of 342 patients, metabolic disease was highly prevalent including nearly 90% overweight.`n`n HS was associated with increased transaminitis and need for intubation, dialysis, and vasopressors.`n`n There was no association between HS and jaundice or portal hypertensive complications.`n`n In a sensitivity analysis including only patients with liver imaging &gt; 30 days before onset of COVID-19, imaging evidence of hepatic steatosis remained associated with disease severity and risk of transaminitis.`n`n`nConclusions`nid=`"Par5`">HS was associated with increased disease severity and transaminitis in COVID-19. HS may be relevant in predicting risk of complications related to COVID-19.`n"
$ws.Range("E3").Value = "[Vincent L.%Chen%vichen@med.umich.edu%0, Fadi%Hawa%NULL%2, Fadi%Hawa%NULL%0, Jeffrey A.%Berinstein%NULL%1, Chanakyaram A.%Reddy%NULL%1, Ihab%Kassab%NULL%1, Kevin D.%Platt%NULL%1, Chia-Yang%Hsu%NULL%1, Calen A.%Steiner%NULL%1, Jeremy%Louissaint%NULL%1, Naresh T.%Gunaratnam%NULL%1, Pratima%Sharma%pratimas@med.umich.edu%2]"
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = "Springer US"
$ws.Range("C4").Value = "Non-alcoholic fatty liver disease (NAFLD) and risk of hospitalization for Covid-19."
$ws.Range("D4").Value = "Background:`nid=`"P2`">Covid-19 disease causes significant morbidity and mortality through increase inflammation and thrombosis.`n Non-alcoholic fatty liver disease and non-alcoholic steatohepatitis are states of chronic inflammation and indicate advanced metabolic disease.`n We sought to understand the risk of hospitalization for Covid-19 associated with NAFLD/NASH.`nMethods:`nid=`"P3`">Retrospective analysis of electronic medical record data of 6,700 adults with a positive SARS-CoV-2 PCR from March 1, 2020 to Aug 25, 2020. Logistic regression and competing risk were used to assess odds of being hospitalized.`n Additional adjustment was added to assess risk of hospitalization among patients with a prescription for metformin use within the 3 months prior to the SARS-CoV-2 PCR result, history of home glucagon-like-peptide 1 receptor agonist (GLP-1 RA) use, and history of metabolic and bariatric surgery (MBS).`n Interactions were assessed by gender and race.`nResults:`nid=`"P4`">A history of NAFLD/NASH was associated with increased odds of admission for Covid-19: logistic regression OR 2.04 (1.55, 2.96, p&lt;0.01), competing risks OR 1.43 (1.09–1.88, p&lt;0.01); and each additional year of having NAFLD/NASH was associated with a significant increased risk of being hospitalized for Covid-19, OR 1.86 (1.43–2.42, p&lt;0.01).`n After controlling for NAFLD/NASH, persons with obesity had decreased odds of hospitalization for Covid-19, OR 0.41 (0.34–0.49, p&lt;0.01).`n NAFLD/NASH increased risk of hospitalization in men and women, and in all racial/ethnic subgroups.`n Mediation treatments for metabolic syndrome were associated with non-significant reduced risk of admission: OR 0.42 (0.18–1.01, p=0.05) for home metformin use and OR 0.40 (0.14–1.17, p=0.10) for home GLP-1RA use.`n MBS was associated with a significant decreased risk of admission: OR 0.22 (0.05–0.98, p&lt;0.05).`nConclusions:`nid=`"P5`">NAFLD/NASH is a significant risk factor for hospitalization for Covid-19, and appears to account for risk attributed to obesity.`n Treatments for metabolic disease mitigated risks from NAFLD/NASH.`n More research is needed to confirm risk associated with visceral adiposity, and patients should be screened for and informed of treatments for metabolic syndrome.`n"
$ws.Range("E4").Value = "[Carolyn T.%Bramante%NULL%0, Christopher J.%Tignanelli%NULL%3, Nirjhar%Dutta%NULL%1, Emma%Jones%NULL%1, Leonardo%Tamaritz%NULL%2, Jeanne%Clark%NULL%1, Genevieve%Melton-Meaux%NULL%1, Michael%Usher%NULL%2, Sayeed%Ikramuddin%NULL%1]"
$ws.Range("F4").Value = "PMC7480063"
$ws.Range("G4").Value = "PMC"
$ws.Range("H4").Value = "1970-01-01"
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = "Cold Spring Harbor Laboratory"
$ws.Range("E5").Value = "[Kenneth I.%Zheng%NULL%0, Feng%Gao%NULL%0, Xiao-Bo%Wang%NULL%3, Qing-Feng%Sun%NULL%2, Ke-Hua%Pan%NULL%2, Ting-Yao%Wang%NULL%2, Hong-Lei%Ma%NULL%2, Yong-Ping%Chen%NULL%3, Wen-Yue%Liu%NULL%2, Jacob%George%NULL%5, Ming-Hua%Zheng%NULL%4]"
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = "Elsevier Inc."
$ws.Range("D6").Value = "Previous studies reported that coronavirus disease 2019 (COVID‐19) was likely to result in liver injury.`n However, few studies investigated liver injury in patients with COVID‐19 with chronic liver diseases.`n We described the clinical features in patients with COVID‐19 with nonalcoholic fatty liver disease (NAFLD).`n Confirmed patients with COVID‐19 from hospitals in 10 cities of Jiangsu Province, China, were retrospectively included between January 18, 2020, and February 26, 2020. The hepatic steatosis index (HSI) was used to defined NAFLD.`n A total of 280 patients with COVID‐19 were enrolled.`n Eighty‐six (30.7%) of 280 patients with COVID‐19 were diagnosed as NAFLD by HSI.`n One hundred (35.7%) patients presented abnormal liver function on admission.`n The median alanine aminotransferase (ALT) levels (34.5 U/L vs.`n 23.0 U/L; P &lt; 0.001) and the proportion of elevated ALT (&gt;40 U/L) (40.7% vs.`n 10.8%; P &lt; 0.001) were significantly higher in patients with NAFLD than in patients without NAFLD on admission.`n The proportion of elevated ALT in patients with NAFLD was also significantly higher than patients without NAFLD (65.1% vs.`n 38.7%; P &lt; 0.001) during hospitalization.`n Multivariate analysis showed that age over 50 years (odds ratio [OR], 2.077; 95% confidence interval [CI], 1.183, 3.648; P = 0.011) and concurrent NAFLD (OR, 2.956; 95% CI, 1.526, 5.726; P = 0.001) were independent risk factors of ALT elevation in patients with COVID‐19, while the atomized inhalation of interferon α‐2b (OR, 0.402; 95% CI, 0.236, 0.683; P = 0.001) was associated with a reduced risk of ALT elevation during hospitalization.`n No patient developed liver failure or death during hospitalization.`n The complications and clinical outcomes were comparable between patients with COVID‐19 with and without NAFLD.`n Conclusion: Patients with NAFLD are more likely to develop liver injury when infected by COVID‐19. However, no patient developed severe liver‐related complications during hospitalization.`n"
$ws.Range("E6").Value = "[Rui%Huang%NULL%0, Li%Zhu%NULL%0, Li%Zhu%NULL%0, Jian%Wang%NULL%5, Leyang%Xue%NULL%1, Longgen%Liu%NULL%2, Xuebing%Yan%NULL%2, Songping%Huang%NULL%1, Yang%Li%NULL%0, Xiaomin%Yan%NULL%1, Biao%Zhang%NULL%1, Tianmin%Xu%NULL%1, Chunyang%Li%NULL%0, Fang%Ji%NULL%2, Fang%Ming%NULL%1, Yun%Zhao%NULL%2, Juan%Cheng%NULL%1, Yinling%Wang%NULL%1, Haiyan%Zhao%NULL%2, Shuqin%Hong%NULL%1, Kang%Chen%NULL%1, Xiang‐an%Zhao%NULL%1, Lei%Zou%NULL%1, Dawen%Sang%NULL%1, Huaping%Shao%NULL%1, Xinying%Guan%NULL%1, Xiaobing%Chen%NULL%1, Yuxin%Chen%NULL%2, Jie%Wei%NULL%1, Chuanwu%Zhu%zhuchw@126.com%3, Chao%Wu%dr.wu@nju.edu.cn%0]"
$ws.Range("I6").Value = ""
$ws.Range("J6").Value = "John Wiley and Sons Inc."
$ws.Range("C7").Value = "Unknown Title"
$ws.Range("D7").Value = "Unknown Abstract"
$ws.Range("E7").Value = "[]"
$ws.Range("F7").Value = "not found"
$ws.Range("G7").Value = "N/A"
$ws.Range("I7").Value = ""
$ws.Range("E8").Value = "[Donghee%Kim%NULL%0, Nia%Adeniji%NULL%1, Nyann%Latt%NULL%1, Sonal%Kumar%NULL%1, Patricia P.%Bloom%NULL%1, Elizabeth S.%Aby%NULL%1, Ponni%Perumalswami%NULL%1, Marina%Roytman%NULL%1, Michael%Li%NULL%1, Alexander S.%Vogel%NULL%1, Andreea M.%Catana%NULL%1, Kara%Wegermann%NULL%1, Rotonya M.%Carr%NULL%1, Costica%Aloman%NULL%1, Vincent L.%Chen%NULL%0, Atoosa%Rabiee%NULL%1, Brett%Sadowski%NULL%1, Veronica%Nguyen%NULL%1, Winston%Dunn%NULL%1, Kenneth D.%Chavin%NULL%1, Kali%Zhou%NULL%1, Blanca%Lizaola-Mayo%NULL%1, Akshata%Moghe%NULL%1, José%Debes%NULL%1, Tzu-Hao%Lee%NULL%1, Andrea D.%Branch%NULL%1, Kathleen%Viveiros%NULL%0, Walter%Chan%NULL%1, David M.%Chascsa%NULL%1, Paul%Kwo%NULL%1, Renumathy%Dhanasekaran%NULL%1]"
$ws.Range("I8").Value = ""
$ws.Range("J8").Value = "by the AGA Institute"
$ws.Range("E9").Value = "[Yu-Jie%Zhou%NULL%0, Kenneth I.%Zheng%NULL%0, Xiao-Bo%Wang%NULL%0, Hua-Dong%Yan%NULL%2, Qing-Feng%Sun%NULL%0, Ke-Hua%Pan%NULL%0, Ting-Yao%Wang%NULL%0, Hong-Lei%Ma%NULL%0, Yong-Ping%Chen%NULL%0, Jacob%George%NULL%0, Ming-Hua%Zheng%NULL%0]"
$ws.Range("I9").Value = ""
$ws.Range("J9").Value = "European Association for the Study of the Liver. Published by Elsevier B.V."
$ws.Range("D10").Value = "The Corona Virus Disease 2019 (COVID‐19) pandemic has attracted increasing worldwide attention.`n While metabolic‐associated fatty liver disease (MAFLD) affects a quarter of world population, its impact on COVID‐19 severity has not been characterized.`n We identified 55 MAFLD patients with COVID‐19, who were 1:1 matched by age, sex and obesity status to non‐aged severe acute respiratory syndrome coronavirus 2 (SARS‐CoV‐2)‐infected patients without MAFLD.`n Our results demonstrate that in patients aged less than 60 years with COVID‐19, MAFLD is associated with an approximately fourfold increase (adjusted odds ratio 4.07, 95% confidence interval 1.20‐13.79, P = .`n02) in the probability for severe disease, after adjusting for confounders.`n Healthcare professionals caring for patients with COVID‐19 need to be aware that there is a positive association between MAFLD and severe illness with COVID‐19."
$ws.Range("E10").Value = "[Yu‐Jie%Zhou%NULL%0, Kenneth I.%Zheng%NULL%0, Kenneth I.%Zheng%NULL%0, Xiao‐Bo%Wang%NULL%4, Xiao‐Bo%Wang%NULL%0, Qing‐Feng%Sun%NULL%3, Ke‐Hua%Pan%NULL%2, Ting‐Yao%Wang%NULL%2, Hong‐Lei%Ma%NULL%1, Yong‐Ping%Chen%NULL%2, Jacob%George%jacob.george@sydney.edu.au%0, Ming‐Hua%Zheng%zhengmh@wmu.edu.cn%2]"
$ws.Range("I10").Value = ""
$ws.Range("J10").Value = "John Wiley and Sons Inc."
$ws.Range("E11").Value = "[Feng%Gao%NULL%0, Kenneth I%Zheng%NULL%1, Xiao‐Bo%Wang%NULL%0, Xiao‐Bo%Wang%NULL%0, Hua‐Dong%Yan%NULL%1, Qing‐Feng%Sun%NULL%0, Ke‐Hua%Pan%NULL%0, Ting‐Yao%Wang%NULL%0, Yong‐Ping%Chen%NULL%0, Jacob%George%jacob.george@sydney.edu.au%0, Ming‐Hua%Zheng%zhengmh@wmu.edu.cn%0]"
$ws.Range("I11").Value = ""
$ws.Range("J11").Value = "John Wiley and Sons Inc."
$ws.Range("E12").Value = "[Ivan%Lopez-Mendez%NULL%0, Jorge%Aquino-Matus%NULL%1, Sofia Murua-Beltrán%Gall%NULL%1, Jose D.%Prieto-Nava%NULL%1, Eva%Juarez-Hernandez%NULL%1, Misael%Uribe%NULL%1, Graciela%Castro-Narro%NULL%1]"
$ws.Range("I12").Value = ""
$ws.Range("J12").Value = "Fundación Clínica Médica Sur, A.C. Published by Elsevier España, S.L.U."
$ws.Range("E13").Value = "[Roberta%Forlano%NULL%0, Benjamin H.%Mullish%NULL%2, Benjamin H.%Mullish%NULL%0, Sujit K.%Mukherjee%NULL%1, Rooshi%Nathwani%NULL%1, Cristopher%Harlow%NULL%1, Peter%Crook%NULL%2, Rebekah%Judge%NULL%2, Anet%Soubieres%NULL%3, Anet%Soubieres%NULL%0, Paul%Middleton%NULL%2, Anna%Daunt%NULL%3, Pablo%Perez-Guzman%NULL%1, Nowlan%Selvapatt%NULL%1, Maud%Lemoine%NULL%1, Ameet%Dhar%NULL%1, Mark R.%Thursz%NULL%1, Shevanthi%Nayagam%NULL%2, Pinelopi%Manousou%NULL%2, Wan-Long%Chuang%NULL%2, Wan-Long%Chuang%NULL%0, NULL%NULL%NULL%0, NULL%NULL%NULL%0]"
$ws.Range("I13").Value = ""
$ws.Range("J13").Value = "Public Library of Science"
$ws.Range("E14").Value = "[Mahmud%Mahamid%NULL%0, Wiliam%Nseir%NULL%1, Tawfik%Khoury%NULL%1, Baker%Mahamid%NULL%1, Adi%Nubania%NULL%1, Kamal%Sub-Laban%NULL%1, Joel%Schifter%NULL%1, Amir%Mari%NULL%1, Wisam%Sbeit%NULL%1, Eran%Goldin%NULL%1]"
$ws.Range("I14").Value = ""
$ws.Range("J14").Value = "Lippincott Williams And Wilkins"
